# Remove oneway-specific bike lane LTS rows.
# The five rows for lanes=3 / oneway=t (B64:B68 = "t") are no longer part
# of the table; delete them entirely so everything below shifts up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A64:G68").EntireRow.Delete()

# Match the author's final view state: scrolled back to the top with
# F11 selected (the row-delete had left the frozen pane/selection
# pointing at cells far below the now-shrunk data range).
$ws.Range("F11").Select() | Out-Null

